$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swaps the content of a single cell (identified by column letter) between
# two rows, using a temporary holding cell far away from the used range.
# Range.Copy(destination) is used instead of Value2 assignment because
# Value2 assignment of a literal "looks-like-a-date" string (e.g.
# "2023-09-21") causes the engine to auto-convert it into a date serial
# number with a new number-format style, which would corrupt the sheet.
# Range.Copy preserves the original cell type/content exactly.
#
# NOTE: this interpreter does not bind PowerShell-style named parameters
# (-col ... -row1 ...) correctly, so plain positional parameters are used
# throughout this script.
function Swap-Cell($col, $row1, $row2, $tempRow) {
    $c1 = $ws.Range($col + $row1)
    $c2 = $ws.Range($col + $row2)
    $tmp = $ws.Range($col + $tempRow)

    # Capture "does a cell exist here" from the ORIGINAL cells before any
    # Copy takes place. Reading .Value2 straight off c1/c2 correctly
    # distinguishes a truly-absent cell (null) from a present-but-empty
    # cell (empty string ""). However, once a value has passed through a
    # Copy() into another range, that distinction is lost (an empty cell
    # copied elsewhere reads back as null too) - so these flags must be
    # captured up front and reused instead of re-inspecting tmp later.
    $c1HasCell = ($c1.Value2 -ne $null)
    $c2HasCell = ($c2.Value2 -ne $null)

    $tmp.ClearContents()

    if ($c1HasCell) {
        $c1.Copy($tmp)
    }

    $c1.ClearContents()

    if ($c2HasCell) {
        $c2.Copy($c1)
    }

    $c2.ClearContents()

    if ($c1HasCell) {
        $tmp.Copy($c2)
    }

    $tmp.ClearContents()
}

function Swap-Rows($row1, $row2, $tempRow, $cols) {
    foreach ($col in $cols) {
        Swap-Cell $col $row1 $row2 $tempRow
    }
}

# Union of columns used by row 39 and row 40 in the original sheet.
$cols3940 = @("A","B","C","D","E","F","G","H","I","J","K","N","P","Q","R","S","T","U","V","W","Y","AA","AD","AE","AF","AG","AT","AW","AX","AY")

# Union of columns used by row 42 and row 43 in the original sheet.
$cols4243 = @("A","B","C","D","E","F","G","H","I","P","Q","R","S","T","U","V","W","Y","AA","AD","AE","AG","AT","AW","AX","AY")

Swap-Rows 39 40 1000 $cols3940
Swap-Rows 42 43 1001 $cols4243
